# Update encryption feature in DG
# Applies the changes described in the target diff to slide 2 of the
# presentation (the "SaveFeatureDiagrams" deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# EMU -> Point helper (PowerPoint COM positions/sizes are in points;
# OOXML stores EMUs, 914400 EMU per inch / 12700 EMU per point). A tiny
# epsilon is added before the division because the host keeps Left/Top/
# Width/Height as single-precision floats internally; without the
# nudge, the point -> EMU round trip on save truncates down and can
# land one EMU short of the intended value.
function EMU($v) { return ($v + 0.4) / 12700.0 }

# ---------------------------------------------------------------------
# 1) "Straight Connector 6" (the short vertical dashed connector near the
#    top of the diagram) grows taller and a bit narrower.
# ---------------------------------------------------------------------
$connector6 = $s.Shapes.Item("Straight Connector 6")
$connector6.Width  = EMU(2664)
$connector6.Height = EMU(3412668)

# ---------------------------------------------------------------------
# 2) Replace the old green dashed "Straight Arrow Connector 48" (and its
#    companion "Close" icon picture, "Graphic 66") with a new blue
#    dashed "Straight Arrow Connector 47" glued to "Rectangle 4".
#    The new connector is produced by duplicating the old one (so it
#    inherits the identical <p:style/> theme refs) and then restyling
#    it, before the original and the icon are deleted.
# ---------------------------------------------------------------------
$oldConn = $s.Shapes.Item("Straight Arrow Connector 48")
$newConn = $oldConn.Duplicate().Item(1)

$newConn.Name = "Straight Arrow Connector 47"
$newConn.HorizontalFlip = $false
$newConn.Left   = EMU(76200)
$newConn.Top    = EMU(5219912)
$newConn.Width  = EMU(1196051)
$newConn.Height = EMU(0)

$newConn.Line.Weight = 1.5
$newConn.Line.DashStyle = 9
$newConn.Line.ForeColor.RGB = 0xC07000
$newConn.Line.BeginArrowheadStyle = 3
$newConn.Line.BeginArrowheadLength = 2
$newConn.Line.BeginArrowheadWidth = 2
$newConn.Line.EndArrowheadStyle = 1
$newConn.Line.EndArrowheadLength = 2
$newConn.Line.EndArrowheadWidth = 2

$rectangle4 = $s.Shapes.Item("Rectangle 4")
$newConn.ConnectorFormat.EndConnect($rectangle4, 2)

$s.Shapes.Item("Graphic 66").Delete()
$oldConn.Delete()

# ---------------------------------------------------------------------
# 3) Nudge "Rectangle 7" to the left.
# ---------------------------------------------------------------------
$rectangle7 = $s.Shapes.Item("Rectangle 7")
$rectangle7.Left = EMU(2949271)
